# "Add files via upload" -- the uploaded workbook gained a new task row for
# team member 邱培松 on the "日期：2018.10.15" plan (row 34, previously blank)
# and the saved cursor position moved from B42 to B40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty row 34 with the new team member / task entry.
$ws.Range("A34").Value = "邱培松"
$ws.Range("B34").Value = "内容:设计数据库E-R模型"

# Match the workbook's saved selection (active cell) at the time of upload.
$ws.Range("B40").Select() | Out-Null
